# Apply the "can_auto_battle" column addition to the Locations sheet.
#
# Summary of the change (reconstructed from the OOXML diff):
#  - A new column is inserted at I ("can_auto_battle"), pushing the former
#    I/J/K columns ("x"/"y"/"type") right to J/K/L.
#  - The new column's header is "can_auto_battle" and, for every data row,
#    its value mirrors the existing "can_players_enter" (column G) value
#    for that row (1 for every location except row 34, which is 0).
#  - The new column gets its own width; the view scroll position / active
#    cell also moved slightly in the source file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new column at I, shifting x/y/type (old I/J/K) to J/K/L.
$ws.Columns("I").Insert()

# --- 2. Header for the freshly inserted column.
$ws.Cells.Item(1, 9).Value = "can_auto_battle"

# --- 3. Populate the new column for every data row (2-35) with the same
#        value as "can_players_enter" (column G) on that row.
$lastRow = 35
for ($r = 2; $r -le $lastRow; $r++) {
    $canEnter = $ws.Cells.Item($r, 7).Value()
    $ws.Cells.Item($r, 9).Value = $canEnter
}

# --- 4. Give the whole new column the same cell style used throughout the
#        rest of the sheet (every other populated cell carries this style).
$ws.Range("I1:I35").Style = $ws.Range("H1").Style

# --- 5. Column width for the new column I.
$ws.Columns("I").ColumnWidth = 17.78

# --- 6. Update the view: active cell moves to G35, and the view scrolls
#        so column D is the left-most visible column.
$ws.Range("G35").Select()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
